# Rename the "prepared" sheet to "icra_results" (this also repoints the
# ICRA2020mostcited defined name, which referred to the "prepared" sheet).
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "icra_results"

# Make icra_results the active/selected tab (was previously icra_bibliography).
$ws1.Activate()
